# Restored from revision #05b0e0d508b8a83d3dfdd86c8a29e69db854e2cb.TEST
# Author: admin. Type: SAVE.
#
# The only substantive content change in this revision is the value of
# cell C10 on the "Rules" sheet, which is restored from 18 back to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
